$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(3, 1, 1.281134, 3.843402, 0.007312702338676299, 0.007312702338676299, 3, 1, 168.1098273333333, 504.329482, 0.2984182258032519, 0.298418225803252, 215.3712155308627, 1938.340939777764, 0.002182243657735073, 0.002182243657735073)
    3 = @(3, 1, 1.281134, 3.843402, 0.007312702338676299, 0.007312702338676299, 3, 1, 163.0062356666667, 489.018707, 0.2893586437755394, 0.2893586437755394, 208.8328307245793, 1879.495476521214, 0.002115993631053589, 0.002115993631053589)
    4 = @(3, 1, 1.281134, 3.843402, 0.007312702338676299, 0.007312702338676299, 3, 1, 165.99353, 497.98059, 0.294661504941043, 0.294661504941043, 212.65995506302, 1913.93959556718, 0.002154771876300243, 0.002154771876300243)
    5 = @(3, 1, 1.281134, 3.843402, 0.007312702338676299, 0.007312702338676299, 3, 1, 66.22673433333334, 198.680203, 0.1175616254801657, 0.1175616254801657, 84.84532106340068, 763.6078895706061, 0.0008596931735873946, 0.0008596931735873946)
    6 = @(3, 1, 164.655884, 493.967652, 0.9398544320918915, 0.9398544320918915, 3, 1, 168.1098273333333, 504.329482, 0.2984182258032519, 0.298418225803252, 27680.27222865736, 249122.4500579163, 0.2804696921381852, 0.2804696921381852)
    7 = @(3, 1, 164.655884, 493.967652, 0.9398544320918915, 0.9398544320918915, 3, 1, 163.0062356666667, 489.018707, 0.2893586437755394, 0.2893586437755394, 26839.93583120733, 241559.422480866, 0.2719550038165395, 0.2719550038165395)
    8 = @(3, 1, 164.655884, 493.967652, 0.9398544320918915, 0.9398544320918915, 3, 1, 165.99353, 497.98059, 0.294661504941043, 0.294661504941043, 27331.81142043052, 245986.3027838747, 0.276938921385706, 0.276938921385706)
    9 = @(3, 1, 164.655884, 493.967652, 0.9398544320918915, 0.9398544320918915, 3, 1, 66.22673433333334, 198.680203, 0.1175616254801657, 0.1175616254801657, 10904.62148608815, 98141.59337479336, 0.1104908147514607, 0.1104908147514607)
    10 = @(3, 1, 9.213772333333333, 27.641317, 0.05259213672418158, 0.05259213672418158, 3, 1, 168.1098273333333, 504.329482, 0.2984182258032519, 0.298418225803252, 1548.92567604531, 13940.33108440779, 0.01569445213243232, 0.01569445213243232)
    11 = @(3, 1, 9.213772333333333, 27.641317, 0.05259213672418158, 0.05259213672418158, 3, 1, 163.0062356666667, 489.018707, 0.2893586437755394, 0.2893586437755394, 1501.902344346347, 13517.12109911712, 0.01521798935576692, 0.01521798935576692)
    12 = @(3, 1, 9.213772333333333, 27.641317, 0.05259213672418158, 0.05259213672418158, 3, 1, 165.99353, 497.98059, 0.294661504941043, 0.294661504941043, 1529.426594226337, 13764.83934803703, 0.01549687815521244, 0.01549687815521244)
    13 = @(3, 1, 9.213772333333333, 27.641317, 0.05259213672418158, 0.05259213672418158, 3, 1, 66.22673433333334, 198.680203, 0.1175616254801657, 0.1175616254801657, 610.1980525274835, 5491.782472747352, 0.006182817080769901, 0.006182817080769901)
    14 = @(1, 0.3333333333333333, 0.042174, 0.126522, 0.0002407288452506406, 0.0002407288452506406, 3, 1, 168.1098273333333, 504.329482, 0.2984182258032519, 0.298418225803252, 7.089863857955999, 63.808774721604, 0.00007183787489936176, 0.00007183787489936178)
    15 = @(1, 0.3333333333333333, 0.042174, 0.126522, 0.0002407288452506406, 0.0002407288452506406, 3, 1, 163.0062356666667, 489.018707, 0.2893586437755394, 0.2893586437755394, 6.874624983005999, 61.871624847054, 0.00006965697217937708, 0.00006965697217937708)
    16 = @(1, 0.3333333333333333, 0.042174, 0.126522, 0.0002407288452506406, 0.0002407288452506406, 3, 1, 165.99353, 497.98059, 0.294661504941043, 0.294661504941043, 7.000611134219999, 63.00550020798, 0.00007093352382427322, 0.00007093352382427322)
    17 = @(1, 0.3333333333333333, 0.042174, 0.126522, 0.0002407288452506406, 0.0002407288452506406, 3, 1, 66.22673433333334, 198.680203, 0.1175616254801657, 0.1175616254801657, 2.793046293774, 25.137416643966, 0.00002830047434762857, 0.00002830047434762857)
}

foreach ($rowNum in $data.Keys) {
    $vals = $data[$rowNum]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $col = 5 + $i
        $ws.Cells.Item($rowNum, $col).Value = $vals[$i]
    }
}